$wb = $excel.ActiveWorkbook

# --- "PN transmission line" sheet: drop the "maximum power flow (MW)" column (E) ---
$wsLine = $wb.Worksheets.Item("PN transmission line")
$wsLine.Range("E1:E42").ClearContents()
$wsLine.Columns("E").ColumnWidth = 8.142857142857142
$wsLine.Range("I2").Select()

# --- "PN bus" sheet: move the selection, no longer the active tab ---
$wsBus = $wb.Worksheets.Item("PN bus")
$wsBus.Range("K22").Select()

# --- "Gen cost in MATPOWER format" sheet: becomes the active/selected tab ---
$wsGen = $wb.Worksheets.Item("Gen cost in MATPOWER format")
$wsGen.Activate()
$wsGen.Range("G12").Select()
